# Insert a new data row at row 195 (pushing existing rows 195:315 down to 196:316)
# and populate it with the new "Apio" price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("195").Insert()

$ws.Range("A195").Value = 3
$ws.Range("B195").Value = "Femacal de La Calera"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44582
$ws.Range("E195").Value = 5
$ws.Range("F195").Value = 100112017
$ws.Range("G195").Value = "Apio"
$ws.Range("H195").Value = "Americana (o)"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 120
$ws.Range("K195").Value = 9000
$ws.Range("L195").Value = 9000
$ws.Range("M195").Value = 9000
$ws.Range("N195").Value = "`$/docena de matas"
$ws.Range("O195").Value = "Pan de Az$([char]0xFA)car"
$ws.Range("P195").Value = 1500
$ws.Range("Q195").Value = 6
$ws.Range("R195").Value = "Hortaliza"
